# Updates the day's Betfair Back/Lay odds sheet:
#  - tweaks a batch of existing odds on rows 2-8
#  - inserts three new Welsh Premiership fixtures as rows 9-11 (pushing the
#    former row 9 "Mexican Liga MX" game down to row 12)
#  - refreshes the (now row 12) Mexican Liga MX odds

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# ---------------------------------------------------------------------------
# Row 2 (Saudi Professional League: Al-Akhdoud vs Al-Kholood Club)
# ---------------------------------------------------------------------------
$ws.Range("J2").Value = 3.25

# ---------------------------------------------------------------------------
# Row 3 (German Bundesliga: Stuttgart vs Eintracht Frankfurt)
# ---------------------------------------------------------------------------
$ws.Range("F3").Value = 1.86
$ws.Range("G3").Value = 1.87
$ws.Range("I3").Value = 4.6
$ws.Range("L3").Value = 1.28
$ws.Range("P3").Value = 2.6
$ws.Range("Q3").Value = 1.59
$ws.Range("W3").Value = 2.14
$ws.Range("X3").Value = 22
$ws.Range("AF3").Value = 13.5
$ws.Range("AO3").Value = 34

# ---------------------------------------------------------------------------
# Row 4 (Saudi Professional League: Dhamk vs Al-Ittihad)
# ---------------------------------------------------------------------------
$ws.Range("G4").Value = 13
$ws.Range("Q4").Value = 1.61

# ---------------------------------------------------------------------------
# Row 5 (Saudi Professional League: Al-Fateh (KSA) vs Al Riyadh SC)
# ---------------------------------------------------------------------------
$ws.Range("U5").Value = 2.2

# ---------------------------------------------------------------------------
# Row 6 (German Bundesliga: Dortmund vs Werder Bremen)
# ---------------------------------------------------------------------------
$ws.Range("F6").Value = 1.37
$ws.Range("L6").Value = 1.22
$ws.Range("R6").Value = 1.77
$ws.Range("T6").Value = 1.76
$ws.Range("V6").Value = 1.11
$ws.Range("W6").Value = 3.6
$ws.Range("X6").Value = 30
$ws.Range("AA6").Value = 300
$ws.Range("AC6").Value = 13.5
$ws.Range("AD6").Value = 34
$ws.Range("AE6").Value = 120
$ws.Range("AF6").Value = 10
$ws.Range("AG6").Value = 10.5
$ws.Range("AI6").Value = 90
$ws.Range("AJ6").Value = 12.5
$ws.Range("AK6").Value = 12.5
$ws.Range("AL6").Value = 27
$ws.Range("AO6").Value = 100

# ---------------------------------------------------------------------------
# Row 7 (German Bundesliga: Hamburger SV vs Leverkusen)
# ---------------------------------------------------------------------------
$ws.Range("I7").Value = 2.12
$ws.Range("J7").Value = 3.9
$ws.Range("K7").Value = 3.95
$ws.Range("L7").Value = 1.29
$ws.Range("V7").Value = 1.89
$ws.Range("W7").Value = 1.36
$ws.Range("X7").Value = 20
$ws.Range("Z7").Value = 15
$ws.Range("AA7").Value = 26
$ws.Range("AB7").Value = 19
$ws.Range("AC7").Value = 8.800000000000001
$ws.Range("AD7").Value = 10.5
$ws.Range("AE7").Value = 19.5
$ws.Range("AF7").Value = 29
$ws.Range("AG7").Value = 15
$ws.Range("AH7").Value = 15.5
$ws.Range("AI7").Value = 29
$ws.Range("AJ7").Value = 70
$ws.Range("AK7").Value = 36
$ws.Range("AL7").Value = 40
$ws.Range("AM7").Value = 65
$ws.Range("AN7").Value = 32
$ws.Range("AO7").Value = 11.5

# ---------------------------------------------------------------------------
# Row 8 (German Bundesliga: Mainz vs FC Heidenheim)
# ---------------------------------------------------------------------------
$ws.Range("J8").Value = 4
$ws.Range("L8").Value = 1.42
$ws.Range("N8").Value = 3.65
$ws.Range("U8").Value = 1.94
$ws.Range("V8").Value = 1.19
$ws.Range("W8").Value = 2.4
$ws.Range("X8").Value = 14
$ws.Range("Z8").Value = 48
$ws.Range("AA8").Value = 190
$ws.Range("AB8").Value = 8
$ws.Range("AC8").Value = 8.6
$ws.Range("AD8").Value = 24
$ws.Range("AE8").Value = 95
$ws.Range("AF8").Value = 9.4
$ws.Range("AG8").Value = 10
$ws.Range("AI8").Value = 95
$ws.Range("AJ8").Value = 16.5
$ws.Range("AK8").Value = 18.5
$ws.Range("AL8").Value = 40
$ws.Range("AM8").Value = 140
$ws.Range("AO8").Value = 120

# ---------------------------------------------------------------------------
# Insert 3 new rows at position 9 - this pushes the existing row 9
# (Mexican Liga MX, Puebla vs Mazatlan FC) down to row 12.
# ---------------------------------------------------------------------------
$ws.Rows("9:11").Insert()

# Columns A..AO in header order, used to translate the data arrays below into
# cell writes.
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO")

# Columns B (Date) and C (Time) hold plain text like "2026-01-13" / "16:45:00"
# in this workbook - force text format so Excel doesn't coerce them into date
# / time serials.
$ws.Range("B9:C11").NumberFormat = "@"

# ---------------------------------------------------------------------------
# New row 9: Welsh Premiership - Caernarfon Town vs Connahs Quay
# ---------------------------------------------------------------------------
$row9 = @("Welsh Premiership","2026-01-13","16:45:00","Caernarfon Town","Connahs Quay",1.02,500,1.02,610,1.02,950,1.01,1.01,1.08,1.23,1.08,1.23,1.08,1.23,1.01,1.01,1.01,1.01,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $colNum = $i + 1
    Set-Cell 9 $colNum $row9[$i]
}

# ---------------------------------------------------------------------------
# New row 10: Welsh Premiership - The New Saints vs Colwyn Bay
# ---------------------------------------------------------------------------
$row10 = @("Welsh Premiership","2026-01-13","16:45:00","The New Saints","Colwyn Bay",1.02,1000,1.02,1000,1.02,950,1.01,1.01,1.08,1.12,1.08,1.12,1.08,1.12,1.01,1.01,1.01,1.01,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $colNum = $i + 1
    Set-Cell 10 $colNum $row10[$i]
}

# ---------------------------------------------------------------------------
# New row 11: Welsh Premiership - Penybont FC vs Cardiff Metropolitan
# ---------------------------------------------------------------------------
$row11 = @("Welsh Premiership","2026-01-13","16:45:00","Penybont FC","Cardiff Metropolitan",1.02,1000,1.02,1000,1.02,950,1.01,1.01,1.08,1.01,1.08,1.01,1.08,1.01,1.01,1.01,1.01,1.01,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $colNum = $i + 1
    Set-Cell 11 $colNum $row11[$i]
}

# ---------------------------------------------------------------------------
# Row 12 (was row 9 before the insert): Mexican Liga MX - Puebla vs Mazatlan FC
# refreshed odds.
# ---------------------------------------------------------------------------
$ws.Range("F12").Value = 1.82
$ws.Range("G12").Value = 2.28
$ws.Range("H12").Value = 3.5
$ws.Range("I12").Value = 5.2
$ws.Range("J12").Value = 3.2
$ws.Range("K12").Value = 950
$ws.Range("P12").Value = 1.87
$ws.Range("Q12").Value = 1.74

Write-Output "Applied Jogos_do_Dia update"
